$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.587.90"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "3.096.31"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'517.37"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'142.60"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "'7.29"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("D11").Value = "'0.374"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "3.620.30"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").Value = "'25.73"
$ws.Range("E14").Value = "  -3.95%  "
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "57.696.43"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "3.089.77"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "'6.13"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("D19").Value = "'13.12"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").Value = "'8.18"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'336.82"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("D24").Value = "'65.78"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("E25").Value = "  +3.67%  "
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").Value = "0.0₃0932"
$ws.Range("E27").Value = "  +3.37%  "
$ws.Range("D28").Value = "'6.44"
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "'20.99"
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("E32").Value = "  -3.06%  "
$ws.Range("D33").Value = "'153.83"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").Value = "'28.03"
$ws.Range("E34").Value = "  +9.89%  "
$ws.Range("D35").Value = "'4.53"
$ws.Range("E35").Value = "  -2.25%  "
$ws.Range("D36").Value = "'5.91"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  -2.55%  "
$ws.Range("D38").Value = "'0.0686"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("D39").Value = "3.134.85"
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("D40").Value = "'36.89"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").Value = "'3.87"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").Value = "'0.672"
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "2.289.81"
$ws.Range("E44").Value = "  +3.74%  "
$ws.Range("E45").Value = "  +2.57%  "
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").Value = "'20.36"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("D48").Value = "'0.950"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").Value = "'5.90"
$ws.Range("E49").Value = "  -3.56%  "
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("E51").Value = "  +1.75%  "
